# Auto-generated Excel COM-interop script to apply the Goblin_Profits market-data refresh.
# The scheduled runner re-pulled Universalis market prices for a handful of leve rows in
# each job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) and recalculated the dependent
# price/profit columns (H:N) for those rows.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 1324.8334
$ws.Range("I43").Value = 1283
$ws.Range("K43").Value = 1283
$ws.Range("M43").Value = -1214

# Row 62
$ws.Range("H62").Value = 9180.923000000001
$ws.Range("I62").Value = 6392.1665
$ws.Range("K62").Value = 6392.1665
$ws.Range("M62").Value = -5768.1665

# Row 65
$ws.Range("H65").Value = 9180.923000000001
$ws.Range("I65").Value = 6392.1665
$ws.Range("K65").Value = 31960.8325
$ws.Range("M65").Value = -28840.8325

# Row 103
$ws.Range("H103").Value = 1362.5
$ws.Range("I103").Value = 1362.5
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 4087.5
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -3501.5
$ws.Range("N103").ClearContents()

# Row 107
$ws.Range("H107").Value = 412236.47
$ws.Range("I107").Value = 463731.62
$ws.Range("J107").Value = 275.33334
$ws.Range("K107").Value = 463731.62
$ws.Range("L107").Value = 275.33334
$ws.Range("M107").Value = -461811.62
$ws.Range("N107").Value = -4115.33334

# Row 137
$ws.Range("H137").Value = 3028.353
$ws.Range("I137").Value = 2752.818
$ws.Range("K137").Value = 8258.454000000002
$ws.Range("M137").Value = -5708.454000000002

# Row 138
$ws.Range("H138").Value = 5278.5537
$ws.Range("I138").Value = 4716.3335
$ws.Range("J138").Value = 5346.02
$ws.Range("K138").Value = 14149.0005
$ws.Range("L138").Value = 16038.06
$ws.Range("M138").Value = -9009.000499999998
$ws.Range("N138").Value = -26318.06

$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Range("H102").Value = 7564.522
$ws.Range("I102").Value = 6143.5
$ws.Range("K102").Value = 6143.5
$ws.Range("M102").Value = -4521.5

# Row 132
$ws.Range("H132").Value = 2444.6365
$ws.Range("I132").Value = 2964.5
$ws.Range("J132").Value = 1534.875
$ws.Range("K132").Value = 8893.5
$ws.Range("L132").Value = 4604.625
$ws.Range("M132").Value = -6363.5
$ws.Range("N132").Value = -9664.625

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 1239
$ws.Range("I22").Value = 1273.2858
$ws.Range("K22").Value = 1273.2858
$ws.Range("M22").Value = -1100.2858

# Row 86
$ws.Range("H86").Value = 43252416
$ws.Range("I86").Value = 2001.2
$ws.Range("J86").Value = 74145570
$ws.Range("K86").Value = 2001.2
$ws.Range("L86").Value = 74145570
$ws.Range("M86").Value = -878.2
$ws.Range("N86").Value = -74147816

# Row 89
$ws.Range("H89").Value = 43252416
$ws.Range("I89").Value = 2001.2
$ws.Range("J89").Value = 74145570
$ws.Range("K89").Value = 10006
$ws.Range("L89").Value = 370727850
$ws.Range("M89").Value = -4390
$ws.Range("N89").Value = -370739082

# Row 94
$ws.Range("H94").Value = 871.06665
$ws.Range("J94").Value = 638.875
$ws.Range("L94").Value = 638.875
$ws.Range("N94").Value = -1540.875

# Row 105
$ws.Range("H105").Value = 3539.8147
$ws.Range("I105").Value = 3290.923
$ws.Range("J105").Value = 10011
$ws.Range("K105").Value = 3290.923
$ws.Range("L105").Value = 10011
$ws.Range("M105").Value = -1543.923
$ws.Range("N105").Value = -13505

# Row 107
$ws.Range("H107").Value = 5668.7144
$ws.Range("I107").Value = 3712.5833
$ws.Range("J107").Value = 8276.888999999999
$ws.Range("K107").Value = 3712.5833
$ws.Range("L107").Value = 8276.888999999999
$ws.Range("M107").Value = -1792.5833
$ws.Range("N107").Value = -12116.889

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5184.625
$ws.Range("I31").Value = 2590
$ws.Range("K31").Value = 2590
$ws.Range("M31").Value = -2295

# Row 34
$ws.Range("H34").Value = 5184.625
$ws.Range("I34").Value = 2590
$ws.Range("K34").Value = 2590
$ws.Range("M34").Value = -2388

# Row 60
$ws.Range("H60").Value = 15748.75
$ws.Range("I60").Value = 1000
$ws.Range("J60").Value = 59995
$ws.Range("K60").Value = 1000
$ws.Range("L60").Value = 59995
$ws.Range("M60").Value = -489
$ws.Range("N60").Value = -61017

# Row 86
$ws.Range("H86").Value = 5999.75
$ws.Range("I86").Value = 6026.846
$ws.Range("K86").Value = 6026.846
$ws.Range("M86").Value = -4903.846

# Row 89
$ws.Range("H89").Value = 5999.75
$ws.Range("I89").Value = 6026.846
$ws.Range("K89").Value = 30134.23
$ws.Range("M89").Value = -24518.23

# Row 135
$ws.Range("H135").Value = 202384.62
$ws.Range("J135").Value = 202384.62
$ws.Range("L135").Value = 202384.62
$ws.Range("N135").Value = -212524.62

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 4626.3335
$ws.Range("I8").Value = 4626.3335
$ws.Range("K8").Value = 13879.0005
$ws.Range("M8").Value = -13740.0005

# Row 37
$ws.Range("H37").Value = 95500
$ws.Range("J37").Value = 95500
$ws.Range("L37").Value = 286500
$ws.Range("N37").Value = -286724

# Row 107
$ws.Range("H107").Value = 1339.1
$ws.Range("J107").Value = 398.53845
$ws.Range("L107").Value = 1195.61535
$ws.Range("N107").Value = -5035.61535

# Row 109
$ws.Range("H109").Value = 10290.75
$ws.Range("I109").Value = 775.3333
$ws.Range("J109").Value = 16000
$ws.Range("K109").Value = 2325.9999
$ws.Range("L109").Value = 48000
$ws.Range("M109").Value = -1285.9999
$ws.Range("N109").Value = -50080

# Row 137
$ws.Range("H137").Value = 5858.05
$ws.Range("I137").Value = 2406.125
$ws.Range("J137").Value = 8159.3335
$ws.Range("K137").Value = 7218.375
$ws.Range("L137").Value = 24478.0005
$ws.Range("M137").Value = -2118.375
$ws.Range("N137").Value = -34678.00049999999

# Row 138
$ws.Range("H138").Value = 3245.5
$ws.Range("J138").Value = 3991
$ws.Range("L138").Value = 11973
$ws.Range("N138").Value = -22253

# Row 139
$ws.Range("H139").Value = 4545.65
$ws.Range("I139").Value = 4434.778
$ws.Range("J139").Value = 4636.364
$ws.Range("K139").Value = 13304.334
$ws.Range("L139").Value = 13909.092
$ws.Range("M139").Value = -8164.334000000001
$ws.Range("N139").Value = -24189.092

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6737.8066
$ws.Range("I70").Value = 6040.4546
$ws.Range("J70").Value = 8442.444
$ws.Range("K70").Value = 6040.4546
$ws.Range("L70").Value = 8442.444
$ws.Range("M70").Value = -5770.4546
$ws.Range("N70").Value = -8982.444

# Row 73
$ws.Range("H73").Value = 6737.8066
$ws.Range("I73").Value = 6040.4546
$ws.Range("J73").Value = 8442.444
$ws.Range("K73").Value = 6040.4546
$ws.Range("L73").Value = 8442.444
$ws.Range("M73").Value = -5104.4546
$ws.Range("N73").Value = -10314.444

# Row 97
$ws.Range("H97").Value = 791.6667
$ws.Range("I97").Value = 619.2857
$ws.Range("K97").Value = 619.2857
$ws.Range("M97").Value = -123.2857

# Row 102
$ws.Range("H102").Value = 8520.272000000001
$ws.Range("I102").Value = 458.57144
$ws.Range("K102").Value = 458.57144
$ws.Range("M102").Value = 1163.42856

# Row 122
$ws.Range("H122").Value = 14481.556
$ws.Range("I122").Value = 17723.334
$ws.Range("K122").Value = 53170.00199999999
$ws.Range("M122").Value = -50720.00199999999

$ws = $wb.Worksheets.Item("LTW")
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("I113").Value = 535.75
$ws.Range("J113").Value = 1055.1538
$ws.Range("K113").Value = 1607.25
$ws.Range("L113").Value = 3165.4614
$ws.Range("M113").Value = 562.75
$ws.Range("N113").Value = -7505.4614

# Row 117
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -59178

# Row 126
$ws.Range("H126").Value = 3027.6875
$ws.Range("I126").Value = 2817.52
$ws.Range("K126").Value = 8452.559999999999
$ws.Range("M126").Value = -5982.559999999999

# Row 136
$ws.Range("H136").Value = 1656.5
$ws.Range("I136").Value = 1249
$ws.Range("K136").Value = 3747
$ws.Range("M136").Value = -1197

